$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 4 has now been assigned to the five backlog items whose Sprint
# column ("C") previously just held the placeholder text "Undecided" /
# "undecided". Replace that placeholder text with the numeric sprint
# number 4 (rows 9-13), leaving the Story Priority column (D) untouched.
$ws.Range("C9:C13").Value = 4

# The active selection in the saved workbook moved to C14.
$ws.Range("C14").Select()
